# Generate Report for Handoff
# Adds two new localization records (06a8ec8f-... and b4e3a99a-...) to every
# sheet of the localization-status workbook:
#   - Overview sheet: new row inserted before the existing 2104fccf row
#     (06a8ec8f) and a new row appended at the end (b4e3a99a).
#   - zh-cn / de-de sheets: same two new rows, with the per-language
#     handoff xliff file name + timestamp columns filled in.
# Existing rows are preserved verbatim; only row positions 5 (inserted),
# 6 (previously row 5) and 7 (new, appended) are touched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview  (columns A:G)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Make room: push the existing "2104fccf" row (row 5) down to row 6, and
# open up a fresh row 7 below it for the "b4e3a99a" row.
$wsOverview.Rows.Item(5).Insert()
$wsOverview.Rows.Item(7).Insert()

# Row 5: 06a8ec8f (new) - same status/date pattern as the 2104fccf record
$wsOverview.Range("A5").Value = "06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d.md"
$wsOverview.Range("B5").Value = "e2e\06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-24 00:39:44"

# Row 6: 2104fccf (existing data, now shifted down one row)
$wsOverview.Range("A6").Value = "2104fccf-d53e-4db5-ac55-e8b2fd450802.md"
$wsOverview.Range("B6").Value = "e2e\2104fccf-d53e-4db5-ac55-e8b2fd450802.md"
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("D6").Value = ""
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-08-24 00:37:55"

# Row 7: b4e3a99a (new, appended)
$wsOverview.Range("A7").Value = "b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6.md"
$wsOverview.Range("B7").Value = "e2e\b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6.md"
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("D7").Value = ""
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-24 00:39:44"

# Rebuild the hyperlinks on column B (Insert() does not carry them along).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48e521a664fa8753bd8e168ce66b1f3f10b8ef1e/e2e/0d966af0-a531-4bee-9c92-fdd766216669.md", [Type]::Missing, [Type]::Missing, "e2e\0d966af0-a531-4bee-9c92-fdd766216669.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b103156c235ecf4810540634a55538bff2e14a0b/e2e/3d4db7e0-e89f-4cf9-b9d7-2ea033e53f2e.md", [Type]::Missing, [Type]::Missing, "e2e\3d4db7e0-e89f-4cf9-b9d7-2ea033e53f2e.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b103156c235ecf4810540634a55538bff2e14a0b/e2e/f58ad3d9-8d5c-4ae9-a4df-0c4550458091.md", [Type]::Missing, [Type]::Missing, "e2e\f58ad3d9-8d5c-4ae9-a4df-0c4550458091.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48e521a664fa8753bd8e168ce66b1f3f10b8ef1e/e2e/06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d.md", [Type]::Missing, [Type]::Missing, "e2e\06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e332398c53e6462c34b18d5133e0ad74d80decd/e2e/2104fccf-d53e-4db5-ac55-e8b2fd450802.md", [Type]::Missing, [Type]::Missing, "e2e\2104fccf-d53e-4db5-ac55-e8b2fd450802.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48e521a664fa8753bd8e168ce66b1f3f10b8ef1e/e2e/b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6.md", [Type]::Missing, [Type]::Missing, "e2e\b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6.md") | Out-Null

# Grow the "Overview" table / AutoFilter to the new A1:G7 extent.
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G7"))

# ---------------------------------------------------------------------
# Shared helper data for the per-language sheets (columns A:P)
# ---------------------------------------------------------------------
function Fill-LangRow($ws, [int]$row, [string]$name, [string]$handoffFile, [string]$handoffDate) {
    $ws.Range("A$row").Value = "$name.md"
    $ws.Range("B$row").Value = ".md"
    $ws.Range("C$row").Value = "Ready for handoff"
    $ws.Range("D$row").Value = "e2e"
    $ws.Range("E$row").Value = "ht"
    $ws.Range("F$row").Value = "False"
    $ws.Range("G$row").Value = $handoffFile
    $ws.Range("H$row").Value = $handoffDate
    $ws.Range("I$row").Value = ""
    $ws.Range("J$row").Value = ""
    $ws.Range("K$row").Value = "0001-01-01 00:00:00"
    $ws.Range("L$row").Value = ""
    $ws.Range("M$row").Value = "True"
    $ws.Range("N$row").Value = ""
    $ws.Range("O$row").Value = "False"
    $ws.Range("P$row").Value = ""
}

# ---------------------------------------------------------------------
# Sheet 2: zh-cn  (columns A:P)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Rows.Item(5).Insert()
$wsZhCn.Rows.Item(7).Insert()

Fill-LangRow $wsZhCn 5 "06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d" "06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d.c9c5cd1eab5f566c21c71fd69737469441a8c50b.zh-cn.xlf" "2016-08-24 00:39:40"
Fill-LangRow $wsZhCn 6 "2104fccf-d53e-4db5-ac55-e8b2fd450802" "2104fccf-d53e-4db5-ac55-e8b2fd450802.f2f2e290a12109c0066bd4c3a6133280d8cf0ab1.zh-cn.xlf" "2016-08-24 00:37:49"
Fill-LangRow $wsZhCn 7 "b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6" "b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6.7599ee1e7d3ed6e853cb1edd12d27f4bc03c22dc.zh-cn.xlf" "2016-08-24 00:39:40"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48e521a664fa8753bd8e168ce66b1f3f10b8ef1e/e2e/0d966af0-a531-4bee-9c92-fdd766216669.md", [Type]::Missing, [Type]::Missing, "0d966af0-a531-4bee-9c92-fdd766216669.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/61b48f377e3bd4d9da6f6488a158282293c86b52/e2e/0d966af0-a531-4bee-9c92-fdd766216669.md", [Type]::Missing, [Type]::Missing, "0d966af0-a531-4bee-9c92-fdd766216669.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b103156c235ecf4810540634a55538bff2e14a0b/e2e/3d4db7e0-e89f-4cf9-b9d7-2ea033e53f2e.md", [Type]::Missing, [Type]::Missing, "3d4db7e0-e89f-4cf9-b9d7-2ea033e53f2e.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b103156c235ecf4810540634a55538bff2e14a0b/e2e/f58ad3d9-8d5c-4ae9-a4df-0c4550458091.md", [Type]::Missing, [Type]::Missing, "f58ad3d9-8d5c-4ae9-a4df-0c4550458091.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48e521a664fa8753bd8e168ce66b1f3f10b8ef1e/e2e/06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d.md", [Type]::Missing, [Type]::Missing, "06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e332398c53e6462c34b18d5133e0ad74d80decd/e2e/2104fccf-d53e-4db5-ac55-e8b2fd450802.md", [Type]::Missing, [Type]::Missing, "2104fccf-d53e-4db5-ac55-e8b2fd450802.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48e521a664fa8753bd8e168ce66b1f3f10b8ef1e/e2e/b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6.md", [Type]::Missing, [Type]::Missing, "b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6.md") | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P7"))

# ---------------------------------------------------------------------
# Sheet 3: de-de  (columns A:P)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Rows.Item(5).Insert()
$wsDeDe.Rows.Item(7).Insert()

Fill-LangRow $wsDeDe 5 "06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d" "06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d.c9c5cd1eab5f566c21c71fd69737469441a8c50b.de-de.xlf" "2016-08-24 00:39:44"
Fill-LangRow $wsDeDe 6 "2104fccf-d53e-4db5-ac55-e8b2fd450802" "2104fccf-d53e-4db5-ac55-e8b2fd450802.f2f2e290a12109c0066bd4c3a6133280d8cf0ab1.de-de.xlf" "2016-08-24 00:37:55"
Fill-LangRow $wsDeDe 7 "b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6" "b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6.7599ee1e7d3ed6e853cb1edd12d27f4bc03c22dc.de-de.xlf" "2016-08-24 00:39:44"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48e521a664fa8753bd8e168ce66b1f3f10b8ef1e/e2e/0d966af0-a531-4bee-9c92-fdd766216669.md", [Type]::Missing, [Type]::Missing, "0d966af0-a531-4bee-9c92-fdd766216669.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6609797bc9decc20f0fc054ead8d6c2598b94241/e2e/0d966af0-a531-4bee-9c92-fdd766216669.md", [Type]::Missing, [Type]::Missing, "0d966af0-a531-4bee-9c92-fdd766216669.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b103156c235ecf4810540634a55538bff2e14a0b/e2e/3d4db7e0-e89f-4cf9-b9d7-2ea033e53f2e.md", [Type]::Missing, [Type]::Missing, "3d4db7e0-e89f-4cf9-b9d7-2ea033e53f2e.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b103156c235ecf4810540634a55538bff2e14a0b/e2e/f58ad3d9-8d5c-4ae9-a4df-0c4550458091.md", [Type]::Missing, [Type]::Missing, "f58ad3d9-8d5c-4ae9-a4df-0c4550458091.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48e521a664fa8753bd8e168ce66b1f3f10b8ef1e/e2e/06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d.md", [Type]::Missing, [Type]::Missing, "06a8ec8f-bb86-4b04-85d7-b06db8f9cc0d.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e332398c53e6462c34b18d5133e0ad74d80decd/e2e/2104fccf-d53e-4db5-ac55-e8b2fd450802.md", [Type]::Missing, [Type]::Missing, "2104fccf-d53e-4db5-ac55-e8b2fd450802.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48e521a664fa8753bd8e168ce66b1f3f10b8ef1e/e2e/b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6.md", [Type]::Missing, [Type]::Missing, "b4e3a99a-ef5f-4a8b-a347-5afd327ad6b6.md") | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P7"))
